{"js": "// Replace the date line and all the three-digit x one-digit multiplication\n// prompts in the practice worksheet with the next day's values.\nconst replacements = [\n  [\"2024-10-10 Thursday\", \"2024-10-11 Friday\"],\n  [\"836\u00d74=\", \"149\u00d79=\"],\n  [\"384\u00d75=\", \"419\u00d73=\"],\n  [\"107\u00d75=\", \"130\u00d78=\"],\n  [\"464\u00d76=\", \"528\u00d72=\"],\n  [\"872\u00d72=\", \"275\u00d72=\"],\n  [\"108\u00d79=\", \"679\u00d78=\"],\n  [\"769\u00d78=\", \"999\u00d74=\"],\n  [\"486\u00d73=\", \"344\u00d76=\"],\n  [\"348\u00d75=\", \"582\u00d78=\"],\n  [\"444\u00d77=\", \"960\u00d75=\"],\n  [\"534\u00d74=\", \"186\u00d72=\"],\n  [\"837\u00d75=\", \"816\u00d79=\"],\n  [\"602\u00d75=\", \"666\u00d77=\"],\n  [\"803\u00d73=\", \"959\u00d79=\"],\n  [\"796\u00d74=\", \"974\u00d79=\"],\n  [\"543\u00d73=\", \"842\u00d76=\"],\n  [\"170\u00d72=\", \"650\u00d76=\"],\n  [\"974\u00d78=\", \"307\u00d76=\"],\n  [\"530\u00d78=\", \"756\u00d75=\"],\n  [\"185\u00d75=\", \"241\u00d74=\"],\n  [\"180\u00d73=\", \"224\u00d73=\"],\n  [\"443\u00d73=\", \"249\u00d79=\"],\n  [\"198\u00d78=\", \"336\u00d73=\"],\n  [\"439\u00d72=\", \"592\u00d74=\"],\n  [\"885\u00d76=\", \"904\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet's date line and all 25 three-digit x one-digit\n# multiplication prompts to the next day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-10 Thursday\", \"2024-10-11 Friday\"),\n    @(\"836\u00d74=\", \"149\u00d79=\"),\n    @(\"384\u00d75=\", \"419\u00d73=\"),\n    @(\"107\u00d75=\", \"130\u00d78=\"),\n    @(\"464\u00d76=\", \"528\u00d72=\"),\n    @(\"872\u00d72=\", \"275\u00d72=\"),\n    @(\"108\u00d79=\", \"679\u00d78=\"),\n    @(\"769\u00d78=\", \"999\u00d74=\"),\n    @(\"486\u00d73=\", \"344\u00d76=\"),\n    @(\"348\u00d75=\", \"582\u00d78=\"),\n    @(\"444\u00d77=\", \"960\u00d75=\"),\n    @(\"534\u00d74=\", \"186\u00d72=\"),\n    @(\"837\u00d75=\", \"816\u00d79=\"),\n    @(\"602\u00d75=\", \"666\u00d77=\"),\n    @(\"803\u00d73=\", \"959\u00d79=\"),\n    @(\"796\u00d74=\", \"974\u00d79=\"),\n    @(\"543\u00d73=\", \"842\u00d76=\"),\n    @(\"170\u00d72=\", \"650\u00d76=\"),\n    @(\"974\u00d78=\", \"307\u00d76=\"),\n    @(\"530\u00d78=\", \"756\u00d75=\"),\n    @(\"185\u00d75=\", \"241\u00d74=\"),\n    @(\"180\u00d73=\", \"224\u00d73=\"),\n    @(\"443\u00d73=\", \"249\u00d79=\"),\n    @(\"198\u00d78=\", \"336\u00d73=\"),\n    @(\"439\u00d72=\", \"592\u00d74=\"),\n    @(\"885\u00d76=\", \"904\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
